$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (axis_flow24) - B2 updated
$ws.Range("B2").Value = 0.1718224919819459

# Row 3 (qbe2) - B3 updated
$ws.Range("B3").Value = 0.1190503615824829

# Row 4 - label changes from nephila25 to axa_adc_assumed, B/C/D updated
$ws.Range("A4").Value = "axa_adc_assumed"
$ws.Range("B4").Value = 0.09569987361192414
$ws.Range("C4").Value = 0.08573758705076666
$ws.Range("D4").Value = -137.3742855

# Row 5 - label changes from axa_adc_assumed to aspen_lpt, B/C/D updated
$ws.Range("A5").Value = "aspen_lpt"
$ws.Range("B5").Value = 0.09478898065925606
$ws.Range("C5").Value = 0.06706731942694223
$ws.Range("D5").Value = -142.5366911

# Row 6 - label changes from aspen_lpt to nephila25, B/C/D updated
$ws.Range("A6").Value = "nephila25"
$ws.Range("B6").Value = 0.07232465211894021
$ws.Range("C6").Value = 0.1574849928203699
$ws.Range("D6").Value = -226.7621411
